# adding validation to program
$wb = $excel.ActiveWorkbook

$wsFileList = $wb.Worksheets.Item("FileList")
$wsConfig   = $wb.Worksheets.Item("Config")

# FileList!A10 : 9 -> 45
$wsFileList.Range("A10").Value = 45

# Config!C2 : "Input/Exhibits" -> "Input/GenericSlip" (new shared string)
$wsConfig.Range("C2").Value = "Input/GenericSlip"

# Update selections / active sheet to match the saved view state:
# FileList is no longer the selected tab, selection moves to A11
$wsFileList.Activate()
$wsFileList.Range("A11").Select()

# Config becomes the active/selected tab, selection moves to B18
$wsConfig.Activate()
$wsConfig.Range("B18").Select()
